$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A180").Value = 178
$ws.Range("B180").Value = "Fiorentina"
$ws.Range("C180").Value = "Napoli"
$ws.Range("D180").Value = 0
$ws.Range("E180").Value = 3
$ws.Range("F180").Value = 1.19
$ws.Range("G180").Value = 1.64
$ws.Range("H180").Value = 1.16
$ws.Range("I180").Value = 1.8
$ws.Range("J180").Value = 0
$ws.Range("K180").Value = 1
$ws.Range("L180").Value = 0.03
$ws.Range("M180").Value = 0.16
$ws.Range("N180").Value = 0.19
$ws.Range("O180").Value = 2

$ws.Range("A181").Value = 179
$ws.Range("B181").Value = "Hellas Verona"
$ws.Range("C181").Value = "Udinese"
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 0
$ws.Range("F181").Value = 1.02
$ws.Range("G181").Value = 0.8100000000000001
$ws.Range("H181").Value = 1.14
$ws.Range("I181").Value = 0.86
$ws.Range("J181").Value = 0
$ws.Range("K181").Value = 0
$ws.Range("L181").Value = 0.12
$ws.Range("M181").Value = 0.05
$ws.Range("N181").Value = 0.16
$ws.Range("O181").Value = 0

$ws.Range("A182").Value = 180
$ws.Range("B182").Value = "Venezia"
$ws.Range("C182").Value = "Empoli"
$ws.Range("D182").Value = 1
$ws.Range("E182").Value = 1
$ws.Range("F182").Value = 0.83
$ws.Range("G182").Value = 1.78
$ws.Range("H182").Value = 0.67
$ws.Range("I182").Value = 1.59
$ws.Range("J182").Value = 0
$ws.Range("K182").Value = 0
$ws.Range("L182").Value = 0.16
$ws.Range("M182").Value = 0.19
$ws.Range("N182").Value = 0.35
$ws.Range("O182").Value = 2

$ws.Range("A183").Value = 181
$ws.Range("B183").Value = "Lecce"
$ws.Range("C183").Value = "Genoa"
$ws.Range("D183").Value = 0
$ws.Range("E183").Value = 0
$ws.Range("F183").Value = 0.5600000000000001
$ws.Range("G183").Value = 1.2
$ws.Range("H183").Value = 0.76
$ws.Range("I183").Value = 0.91
$ws.Range("J183").Value = 0
$ws.Range("K183").Value = 1
$ws.Range("L183").Value = 0.2
$ws.Range("M183").Value = 0.29
$ws.Range("N183").Value = 0.49
$ws.Range("O183").Value = 1

$ws.Range("A184").Value = 182
$ws.Range("B184").Value = "Monza"
$ws.Range("C184").Value = "Cagliari"
$ws.Range("D184").Value = 1
$ws.Range("E184").Value = 2
$ws.Range("F184").Value = 1.58
$ws.Range("G184").Value = 0.85
$ws.Range("H184").Value = 1.36
$ws.Range("I184").Value = 1.09
$ws.Range("J184").Value = 1
$ws.Range("K184").Value = 0
$ws.Range("L184").Value = 0.22
$ws.Range("M184").Value = 0.24
$ws.Range("N184").Value = 0.46
$ws.Range("O184").Value = 2

$ws.Range("A185").Value = 183
$ws.Range("B185").Value = "Roma"
$ws.Range("C185").Value = "Lazio"
$ws.Range("D185").Value = 2
$ws.Range("E185").Value = 0
$ws.Range("F185").Value = 1.18
$ws.Range("G185").Value = 0.87
$ws.Range("H185").Value = 1.04
$ws.Range("I185").Value = 0.87
$ws.Range("J185").Value = 0
$ws.Range("K185").Value = 0
$ws.Range("L185").Value = 0.14
$ws.Range("M185").Value = 0
$ws.Range("N185").Value = 0.14
$ws.Range("O185").Value = 2

$ws.Range("A186").Value = 184
$ws.Range("B186").Value = "Torino"
$ws.Range("C186").Value = "Parma"
$ws.Range("D186").Value = 0
$ws.Range("E186").Value = 0
$ws.Range("F186").Value = 1.08
$ws.Range("G186").Value = 0.38
$ws.Range("H186").Value = 1.6
$ws.Range("I186").Value = 0.42
$ws.Range("J186").Value = 0
$ws.Range("K186").Value = 0
$ws.Range("L186").Value = 0.52
$ws.Range("M186").Value = 0.04
$ws.Range("N186").Value = 0.5600000000000001
$ws.Range("O186").Value = 0

$ws.Range("A187").Value = 185
$ws.Range("B187").Value = "Lazio"
$ws.Range("C187").Value = "Como"
$ws.Range("D187").Value = 1
$ws.Range("E187").Value = 1
$ws.Range("F187").Value = 0.53
$ws.Range("G187").Value = 1.98
$ws.Range("H187").Value = 0.95
$ws.Range("I187").Value = 1.39
$ws.Range("J187").Value = 0
$ws.Range("K187").Value = 0
$ws.Range("L187").Value = 0.42
$ws.Range("M187").Value = 0.59
$ws.Range("N187").Value = 1.01
$ws.Range("O187").Value = 2

$ws.Range("A188").Value = 186
$ws.Range("B188").Value = "Empoli"
$ws.Range("C188").Value = "Lecce"
$ws.Range("D188").Value = 1
$ws.Range("E188").Value = 3
$ws.Range("F188").Value = 1.43
$ws.Range("G188").Value = 1.46
$ws.Range("H188").Value = 1.61
$ws.Range("I188").Value = 1.24
$ws.Range("J188").Value = 0
$ws.Range("K188").Value = 0
$ws.Range("L188").Value = 0.18
$ws.Range("M188").Value = 0.22
$ws.Range("N188").Value = 0.4
$ws.Range("O188").Value = 4

$ws.Range("A189").Value = 187
$ws.Range("B189").Value = "Milan"
$ws.Range("C189").Value = "Cagliari"
$ws.Range("D189").Value = 1
$ws.Range("E189").Value = 1
$ws.Range("F189").Value = 2.85
$ws.Range("G189").Value = 0.43
$ws.Range("H189").Value = 2.77
$ws.Range("I189").Value = 0.45
$ws.Range("J189").Value = 1
$ws.Range("K189").Value = 0
$ws.Range("L189").Value = 0.08
$ws.Range("M189").Value = 0.02
$ws.Range("N189").Value = 0.1
$ws.Range("O189").Value = 1

$ws.Range("A190").Value = 188
$ws.Range("B190").Value = "Torino"
$ws.Range("C190").Value = "Juventus"
$ws.Range("D190").Value = 1
$ws.Range("E190").Value = 1
$ws.Range("F190").Value = 0.59
$ws.Range("G190").Value = 0.9
$ws.Range("H190").Value = 0.73
$ws.Range("I190").Value = 1.31
$ws.Range("J190").Value = 0
$ws.Range("K190").Value = 0
$ws.Range("L190").Value = 0.14
$ws.Range("M190").Value = 0.41
$ws.Range("N190").Value = 0.55
$ws.Range("O190").Value = 2

$ws.Range("A191").Value = 189
$ws.Range("B191").Value = "Udinese"
$ws.Range("C191").Value = "Atalanta"
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 0
$ws.Range("F191").Value = 1.83
$ws.Range("G191").Value = 0.23
$ws.Range("H191").Value = 1.61
$ws.Range("I191").Value = 0.3
$ws.Range("J191").Value = 1
$ws.Range("K191").Value = 0
$ws.Range("L191").Value = 0.22
$ws.Range("M191").Value = 0.07000000000000001
$ws.Range("N191").Value = 0.29
$ws.Range("O191").Value = 1

$ws.Range("A192").Value = 190
$ws.Range("B192").Value = "Bologna"
$ws.Range("C192").Value = "Roma"
$ws.Range("D192").Value = 2
$ws.Range("E192").Value = 2
$ws.Range("F192").Value = 2.06
$ws.Range("G192").Value = 1.4
$ws.Range("H192").Value = 2.05
$ws.Range("I192").Value = 1.45
$ws.Range("J192").Value = 1
$ws.Range("K192").Value = 1
$ws.Range("L192").Value = 0.01
$ws.Range("M192").Value = 0.05
$ws.Range("N192").Value = 0.06
$ws.Range("O192").Value = 2

$ws.Range("A193").Value = 191
$ws.Range("B193").Value = "Genoa"
$ws.Range("C193").Value = "Parma"
$ws.Range("D193").Value = 1
$ws.Range("E193").Value = 0
$ws.Range("F193").Value = 1.19
$ws.Range("G193").Value = 0.46
$ws.Range("H193").Value = 1.04
$ws.Range("I193").Value = 0.7
$ws.Range("J193").Value = 0
$ws.Range("K193").Value = 0
$ws.Range("L193").Value = 0.15
$ws.Range("M193").Value = 0.24
$ws.Range("N193").Value = 0.38
$ws.Range("O193").Value = 1

$ws.Range("A194").Value = 192
$ws.Range("B194").Value = "Napoli"
$ws.Range("C194").Value = "Hellas Verona"
$ws.Range("D194").Value = 2
$ws.Range("E194").Value = 0
$ws.Range("F194").Value = 1.18
$ws.Range("G194").Value = 0.27
$ws.Range("H194").Value = 1.51
$ws.Range("I194").Value = 0.31
$ws.Range("J194").Value = 0
$ws.Range("K194").Value = 0
$ws.Range("L194").Value = 0.33
$ws.Range("M194").Value = 0.04
$ws.Range("N194").Value = 0.38
$ws.Range("O194").Value = 2

$ws.Range("A195").Value = 193
$ws.Range("B195").Value = "Venezia"
$ws.Range("C195").Value = "Inter"
$ws.Range("D195").Value = 0
$ws.Range("E195").Value = 1
$ws.Range("F195").Value = 0.41
$ws.Range("G195").Value = 1.93
$ws.Range("H195").Value = 0.57
$ws.Range("I195").Value = 2.03
$ws.Range("J195").Value = 0
$ws.Range("K195").Value = 0
$ws.Range("L195").Value = 0.16
$ws.Range("M195").Value = 0.1
$ws.Range("N195").Value = 0.27
$ws.Range("O195").Value = 1

# Copy the bold/bordered/centered style from the last existing data row's
# column-A cell (A179) down across the newly appended rows (A180:A195),
# matching the existing pattern where only column A carries style index 1.
$ws.Range("A179").Copy()
$ws.Range("A180:A195").PasteSpecial(-4122)
$excel.CutCopyMode = 0
